# This workbook lists species observations ("artfynd"). The edit reorders
# the four data rows (rows 3-6): row 3 swaps places with row 6, and row 4
# swaps places with row 5. Only the columns whose values actually differ
# between the two rows being swapped are written, so columns that already
# hold identical data in both rows (and could otherwise be re-interpreted,
# e.g. date-looking text turned into real dates) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns present in the data rows that can hold differing values.
$columns = @("A","B","C","D","E","F","G","H","I","P","Q","R","S","T","U","V","W","Z","AB","AD","AE","AG","AT","AW","AX","AY")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $columns) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")
        $valA = $cellA.Value2
        $valB = $cellB.Value2

        # Only touch cells whose values actually differ, to avoid Excel
        # re-interpreting re-written text (e.g. dates) unnecessarily.
        $same = $false
        if ($valA -eq $null -and $valB -eq $null) {
            $same = $true
        } elseif ($valA -ne $null -and $valB -ne $null -and $valA -eq $valB) {
            $same = $true
        }

        if (-not $same) {
            $cellA.Value2 = $valB
            $cellB.Value2 = $valA
        }
    }
}

Swap-Rows 3 6
Swap-Rows 4 5
